# AHDT1_AHD_VAR.docx edit script
# Applies: title merge, body text merges (run consolidation, no visible text
# change), the Test Limitations detection-limit rewrite (2% -> 4%, ASXL1/CEBPA
# exception -> JAK2 exception), and the report date update.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
    return $ok
}

# 1. Title: merge "HAEMATOLOGICAL MALIGNANCY " + "GENE PANEL REPORT" runs
Replace-Text "HAEMATOLOGICAL MALIGNANCY GENE PANEL REPORT" "HAEMATOLOGICAL MALIGNANCY GENE PANEL REPORT"

# 2. "Somatic variant analysis of 80 genes..." - merge trailing runs
Replace-Text " genes with clinical significance in haematological malignancy plus analysis of potential germline variants in the DDX41 gene. Refer to Panel Summary for gene list." " genes with clinical significance in haematological malignancy plus analysis of potential germline variants in the DDX41 gene. Refer to Panel Summary for gene list."

# 3. "Please note, variant origin..." merge
Replace-Text "Please note, variant origin (somatic or germline) cannot be determined by this assay. Variant origin is assumed here based on ancillary information (e.g. population databases, literature, variant read frequency) for the purpose of clinical interpretation however testing of a germline sample may be recommended in some circumstances." "Please note, variant origin (somatic or germline) cannot be determined by this assay. Variant origin is assumed here based on ancillary information (e.g. population databases, literature, variant read frequency) for the purpose of clinical interpretation however testing of a germline sample may be recommended in some circumstances."

# 4. "VRF - variant read frequency" merge
Replace-Text "VRF – variant read frequency " "VRF – variant read frequency "

# 5-8. Test Methodology paragraph merges (AllHaem / NovaSeq / Seqliner-Nextflow / PathOS / gnomAD / ClinVar)
Replace-Text "panel (Peter MacCallum Cancer Centre AllHaem DNA Twist v1, design ID" "panel (Peter MacCallum Cancer Centre AllHaem DNA Twist v1, design ID"
Replace-Text ") and sequenced on an Illumina NovaSeq 6000 with 150 bp paired end reads." ") and sequenced on an Illumina NovaSeq 6000 with 150 bp paired end reads."
Replace-Text "A custom Seqliner/Nextflow-based analysis pipeline is used to generate aligned reads and call variants (single nucleotide variants and short insertions or deletions) against the hg19 human reference genome." "A custom Seqliner/Nextflow-based analysis pipeline is used to generate aligned reads and call variants (single nucleotide variants and short insertions or deletions) against the hg19 human reference genome."
Replace-Text "Variants are analysed using PathOS software (Peter Mac) and described according to HGVS nomenclature version 19.01 (http://varnomen.hgvs.org/) with minor differences in accordance with Peter MacCallum Cancer Centre Molecular Pathology departmental policy. The following population variation and cancer or genetic disease databases are commonly used in addition to literature review to assist with variant interpretation: the Genome Aggregation Database (gnomAD; gnomad.broadinstitute.org), the Catalogue of Somatic Mutations in Cancer (COSMIC; cancer.sanger.ac.uk), ClinVar (ncbi.nlm.nih.gov/clinvar) and the IARC TP53 Database (p53.iarc.fr). " "Variants are analysed using PathOS software (Peter Mac) and described according to HGVS nomenclature version 19.01 (http://varnomen.hgvs.org/) with minor differences in accordance with Peter MacCallum Cancer Centre Molecular Pathology departmental policy. The following population variation and cancer or genetic disease databases are commonly used in addition to literature review to assist with variant interpretation: the Genome Aggregation Database (gnomAD; gnomad.broadinstitute.org), the Catalogue of Somatic Mutations in Cancer (COSMIC; cancer.sanger.ac.uk), ClinVar (ncbi.nlm.nih.gov/clinvar) and the IARC TP53 Database (p53.iarc.fr). "

# 9. "Variant origin (i.e. somatic or germline)..." merge
Replace-Text "Variant origin (i.e. somatic or germline) is assumed based on ancillary information (e.g. population databases, literature, variant read frequency) for the purpose of clinical interpretation. All assumed somatic variants are reported (and generally considered clinically significant). Variants of uncertain origin are also reported, as are likely benign germline polymorphisms if sufficiently rare and otherwise undescribed. Testing of a non-haematological specimen may be recommended to evaluate variant origin. Recurrent population variants are not reported." "Variant origin (i.e. somatic or germline) is assumed based on ancillary information (e.g. population databases, literature, variant read frequency) for the purpose of clinical interpretation. All assumed somatic variants are reported (and generally considered clinically significant). Variants of uncertain origin are also reported, as are likely benign germline polymorphisms if sufficiently rare and otherwise undescribed. Testing of a non-haematological specimen may be recommended to evaluate variant origin. Recurrent population variants are not reported."

# 10. "The detection limit for FLT3-ITDs is approximately 1%." merge
Replace-Text "The detection limit for FLT3-ITDs is approximately 1%. " "The detection limit for FLT3-ITDs is approximately 1%. "

# 11. DIAGNOSTIC category merge
Replace-Text " (the variant either defines a diagnostic category or is sufficiently specific for the clinical context to contribute to diagnostic subcategorisation), " " (the variant either defines a diagnostic category or is sufficiently specific for the clinical context to contribute to diagnostic subcategorisation), "

# 12. PROGNOSTIC category merge
Replace-Text " (the variant has been associated in large trials/series with inferior or superior outcomes in either the context of a specific therapy or independent of therapy. Note this does not take into account interaction between prognostic variants present in the individual patient. Relevant pairwise interactions are presented in the clinical summary), " " (the variant has been associated in large trials/series with inferior or superior outcomes in either the context of a specific therapy or independent of therapy. Note this does not take into account interaction between prognostic variants present in the individual patient. Relevant pairwise interactions are presented in the clinical summary), "

# 13. DRUG TARGET category merge
Replace-Text " (the variant or variant class is specifically targeted by a therapeutic agent, this category only includes therapeutic agents that are clinically advanced and generally available through either reimbursement or clinical trials [i.e. not early stage investigational agents]), " " (the variant or variant class is specifically targeted by a therapeutic agent, this category only includes therapeutic agents that are clinically advanced and generally available through either reimbursement or clinical trials [i.e. not early stage investigational agents]), "

# 14. DRUG RESISTANCE category merge
Replace-Text " (the variant is specifically associated with resistance to a targeted agent [i.e. does not include non-specific resistance to non-targeted therapies]), " " (the variant is specifically associated with resistance to a targeted agent [i.e. does not include non-specific resistance to non-targeted therapies]), "

# 15. MRD MARKER category merge
Replace-Text " (the variant is an established biomarker for which assessment at MRD sensitivity after therapy is accepted practice). If the variant is not categorised into any of the above categories it is assigned " " (the variant is an established biomarker for which assessment at MRD sensitivity after therapy is accepted practice). If the variant is not categorised into any of the above categories it is assigned "

# 16. "Germline variant analysis" merge
Replace-Text "Germline variant analysis" "Germline variant analysis"

# 17. class 3/4/5 merge
Replace-Text " with class 3 (uncertain significance), class 4 (likely pathogenic) and class 5 (pathogenic) variants reported only. Please note however that germline confirmation is required for all potential clinically significant DDX41 variants." " with class 3 (uncertain significance), class 4 (likely pathogenic) and class 5 (pathogenic) variants reported only. Please note however that germline confirmation is required for all potential clinically significant DDX41 variants."

# 18. BIG CHANGE: detection limit paragraph rewrite
$oldDetection = "The detection limit of this assay for specimens sequenced to the target read depth of 250x is a variant allele frequency (VAF) of approximately 2% with the exception of ASXL1 c.1934dup;p.Gly646Trpfs*12 (detection limit ~ 5%-10%) and CEBPA (detection limit ~ 10%). This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length), including FLT3-ITDs and UBTF-TDs, are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. This assay does not distinguish between somatic and germline variants. In addition, the clonal origin of somatic variants (i.e. disease compartment or cell lineage) cannot be determined. "
$newDetection = "The detection limit of this assay for specimens sequenced to the target read depth of 250x is a variant allele frequency (VAF) of approximately 4% with the exception of JAK2 c.1849G>T;p.(Val617Phe) (detection limit ~ 1%). This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length), including FLT3-ITDs and UBTF-TDs, are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. This assay does not distinguish between somatic and germline variants. In addition, the clonal origin of somatic variants (i.e. disease compartment or cell lineage) cannot be determined. "
Replace-Text $oldDetection $newDetection

# 19. germline zygosity paragraph merge
Replace-Text ", variant zygosity is assumed to be either heterozygous or homozygous in the germline based on allele frequency for the purpose of clinical interpretation however, the possibilities of hemizygosity or somatic acquisition are not excluded. In haematological specimens, the possibility of a false negative germline result due to loss of the mutant allele through a somatic reversion event cannot be excluded. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient." ", variant zygosity is assumed to be either heterozygous or homozygous in the germline based on allele frequency for the purpose of clinical interpretation however, the possibilities of hemizygosity or somatic acquisition are not excluded. In haematological specimens, the possibility of a false negative germline result due to loss of the mutant allele through a somatic reversion event cannot be excluded. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient."

# 20. "A separate assay may have been performed..." merge
Replace-Text "detected with this assay. A separate assay may have been performed, result included in Test Results if sample tested." "detected with this assay. A separate assay may have been performed, result included in Test Results if sample tested."

# 21. "Please note variants may not be optimally detected..." merge
Replace-Text "Please note variants may not be optimally detected in genes with less than 100% coverage. The gene coverage above is considered acceptable given the available information about the clinical context, however please contact the laboratory for further advice should specific genes covered at less than 100% require full coverage. A list of regions with suboptimal coverage is available upon request." "Please note variants may not be optimally detected in genes with less than 100% coverage. The gene coverage above is considered acceptable given the available information about the clinical context, however please contact the laboratory for further advice should specific genes covered at less than 100% require full coverage. A list of regions with suboptimal coverage is available upon request."

# 22. "Please contact the laboratory on 03 8559 7284..." merge
Replace-Text "Please contact the laboratory on 03 8559 7284 if you wish to discuss this report further." "Please contact the laboratory on 03 8559 7284 if you wish to discuss this report further."

# 23. Report date field result update
Replace-Text "16-Sep-2024" "7-Oct-2024"

Write-Output "Done"
